# Fruta / hortaliza, semanal
# Insert a new weekly record as row 98, pushing the existing rows 98-104
# down to 99-105 (row 105 becomes an exact copy of the former row 104).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 98 (shifts rows 98:104 down to 99:105,
# Excel copies formatting such as the date style from the row above).
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new weekly entry.
$ws.Range("A98").Value = 7
$ws.Range("B98").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C98").Value = "Ñuble"
$ws.Range("D98").Value = 44918
$ws.Range("E98").Value = 16
$ws.Range("F98").Value = 100112030
$ws.Range("G98").Value = "Poroto granado"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 60
$ws.Range("K98").Value = 35000
$ws.Range("L98").Value = 35000
$ws.Range("M98").Value = 35000
$ws.Range("N98").Value = "$/saco 25 kilos"
$ws.Range("O98").Value = "Región del Maule"
$ws.Range("P98").Value = 1400
$ws.Range("Q98").Value = 25
$ws.Range("R98").Value = "Hortaliza"
